$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 (shifts shaik etc down) for "sathiya" (id 303)
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 2).Value = 303
$ws.Cells.Item(2, 3).Value = "sathiya"
$ws.Cells.Item(2, 4).Value = "sathiya@gmail.com"
$ws.Cells.Item(2, 5).Value = "globalTiger"
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = "55667788"
$ws.Cells.Item(2, 7).Value = "DM_rejected"

# Insert a new row above current row 4 (after shaik, before aman) for "saisameer" (id 306)
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 2).Value = 306
$ws.Cells.Item(4, 3).Value = "saisameer"
$ws.Cells.Item(4, 4).Value = "saidameer@gmail.com"
$ws.Cells.Item(4, 5).Value = "tcs"
$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = "123456789"
$ws.Cells.Item(4, 7).Value = "DM_rejected"

# Append new row 7 for "cleveland1" (id 318)
$ws.Cells.Item(7, 2).Value = 318
$ws.Cells.Item(7, 3).Value = "cleveland1"
$ws.Cells.Item(7, 4).Value = "cleveland1@gmail.com"
$ws.Cells.Item(7, 5).Value = "nasa corporation"
$ws.Cells.Item(7, 6).NumberFormat = "@"
$ws.Cells.Item(7, 6).Value = "55667788"
$ws.Cells.Item(7, 7).Value = "DM_TBS"

# Append new row 8 for "surmak" (id 319)
$ws.Cells.Item(8, 2).Value = 319
$ws.Cells.Item(8, 3).Value = "surmak"
$ws.Cells.Item(8, 4).Value = "surmka@gmail.com"
$ws.Cells.Item(8, 5).Value = "umbrala corporation"
$ws.Cells.Item(8, 6).NumberFormat = "@"
$ws.Cells.Item(8, 6).Value = "55667788"
$ws.Cells.Item(8, 7).Value = "DM_scheduled"
